# Swap the species-observation data between row 4 and row 6, leaving the
# shared/location columns (C, D, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "I", "J", "Q", "R")

foreach ($col in $cols) {
    $cell4 = $ws.Range($col + "4")
    $cell6 = $ws.Range($col + "6")

    $v4 = $cell4.Value2
    $v6 = $cell6.Value2

    $cell4.Value = $v6
    $cell6.Value = $v4
}
